$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Patient identification (row 6)
$ws.Range("A6").Value = "HIPOLITO"
$ws.Range("C6").Value = "SAMUEL "
$ws.Range("E6").Value = "JUÀREZ"
$ws.Range("G6").Value = "DEL AGUILA "
$ws.Range("I6").Value = "2017-41934/201761796"

# Fecha de nacimiento / Edad (row 12)
$ws.Range("A12").Value = "1970-08-16"
$ws.Range("F12").Value = "47"

# Ocupación / Nacionalidad / No. de Cédula (row 14)
$ws.Range("D14").Value = ""
$ws.Range("F14").Value = "GUATEMALTECO"
$ws.Range("H14").Value = ""

# Emergency contact (row 20)
$ws.Range("A20").Value = "LORENA GOMEZ "
$ws.Range("F20").Value = "ESPOSA "
$ws.Range("H20").Value = "El mosquito z3, San Marcos"
$ws.Range("J20").Value = ""
